$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.205.06"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.834.94"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +1.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.70"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4715"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3687"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07427"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8825"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.43"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "1.822.27"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07338"
$ws.Range("E13").Value = "  +3.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.484"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.92"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.571"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.012"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008793"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "27.227.74"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.69"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("D24").Value = "2.065.06"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.905"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.56"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.62"
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.162"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.282"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.83"
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08933"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.172"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.551"
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.944"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.011"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.103"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05344"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01962"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.008"
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5350"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.545"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4952"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.52"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.011"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.99"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.671"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06330"
$ws.Range("E51").Value = "  +0.72%  "

# Swap rows 41 and 42 (FraxShare <-> RenderToken)
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.428"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.351"
$ws.Range("E42").Value = "  +1.22%  "
